$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.906.67"
$ws.Range("E2").Value = "  -2.13%  "

$ws.Range("D3").Value = "'3.289.43"
$ws.Range("E3").Value = "  -0.70%  "

$ws.Range("D5").Value = "'573.49"
$ws.Range("E5").Value = "  -0.77%  "

$ws.Range("D6").Value = "'178.15"
$ws.Range("E6").Value = "  -4.42%  "

$ws.Range("D7").Value = "'0.631"
$ws.Range("E7").Value = "  +4.74%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  -2.53%  "

$ws.Range("D10").Value = "'6.70"
$ws.Range("E10").Value = "  +0.85%  "

$ws.Range("E11").Value = "  -2.45%  "

$ws.Range("D12").Value = "'3.862.59"
$ws.Range("E12").Value = "  -0.56%  "

$ws.Range("E13").Value = "  -3.50%  "

$ws.Range("D14").Value = "'26.65"
$ws.Range("E14").Value = "  -2.96%  "

$ws.Range("D15").Value = "'66.007.60"
$ws.Range("E15").Value = "  -2.31%  "

$ws.Range("E16").Value = "  -2.29%  "

$ws.Range("D17").Value = "'3.274.89"
$ws.Range("E17").Value = "  -1.30%  "

$ws.Range("D18").Value = "'437.65"
$ws.Range("E18").Value = "  -1.31%  "

$ws.Range("E19").Value = "  -2.34%  "

$ws.Range("D20").Value = "'13.28"
$ws.Range("E20").Value = "  -1.90%  "

$ws.Range("E21").Value = "  -4.66%  "

$ws.Range("D22").Value = "'72.52"
$ws.Range("E22").Value = "  -1.85%  "

$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.23%  "

$ws.Range("B24").Value = "WrappedeETH"
$ws.Range("C24").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D24").Value = "'3.434.59"
$ws.Range("E24").Value = "  -0.63%  "

$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").Value = "'0.511"
$ws.Range("E25").Value = "  -0.41%  "

$ws.Range("E26").Value = "  -4.95%  "

$ws.Range("E27").Value = "  +2.92%  "

$ws.Range("D28").Value = "'8.93"
$ws.Range("E28").Value = "  -1.30%  "

$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("E30").Value = "  -2.03%  "

$ws.Range("D31").Value = "'22.38"
$ws.Range("E31").Value = "  -2.25%  "

$ws.Range("E32").Value = "  +0.12%  "

$ws.Range("D33").Value = "'5.16"
$ws.Range("E33").Value = "  -3.52%  "

$ws.Range("D34").Value = "'6.63"
$ws.Range("E34").Value = "  -2.59%  "

$ws.Range("E35").Value = "  -3.97%  "

$ws.Range("D36").Value = "'158.79"
$ws.Range("E36").Value = "  -2.55%  "

$ws.Range("E37").Value = "  -4.45%  "

$ws.Range("E38").Value = "  -1.66%  "

$ws.Range("E39").Value = "  -3.60%  "

$ws.Range("D40").Value = "'2.778.83"
$ws.Range("E40").Value = "  -0.32%  "

$ws.Range("E41").Value = "  -1.20%  "

$ws.Range("E42").Value = "  -2.57%  "

$ws.Range("D43").Value = "'40.36"
$ws.Range("E43").Value = "  +0.68%  "

$ws.Range("D44").Value = "'6.04"
$ws.Range("E44").Value = "  -3.29%  "

$ws.Range("E45").Value = "  -2.22%  "

$ws.Range("D46").Value = "'2.30"
$ws.Range("E46").Value = "  -4.66%  "

$ws.Range("D47").Value = "'321.49"
$ws.Range("E47").Value = "  -2.29%  "

$ws.Range("D48").Value = "'23.49"
$ws.Range("E48").Value = "  -5.64%  "

$ws.Range("D49").Value = "'0.0269"
$ws.Range("E49").Value = "  -1.99%  "

$ws.Range("D50").Value = "'0.102"
$ws.Range("E50").Value = "  +2.77%  "

$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  +0.10%  "
